# UC007 - Listar Autorizações de Pagamento Pendentes
# Version bump 1.0 -> 1.2.5, typo/accent/punctuation fixes, and
# re-sequencing of TC2/TC3/TC4 step content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header: Version ---
$ws.Range("D2").Value = "1.2.5"

# --- TC1 (rows 6-12) ---
$ws.Range("B8").Value  = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("D10").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B11").Value = "Chefe Seleciona uma diária apta para pagamento."
$ws.Range("D11").Value = "SYSTEM Destaca a diária selecionada."

# --- TC2 (rows 15-20) ---
$ws.Range("B17").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B19").Value = "Chefe Clica para exibir a lista de diárias (solicitações) aptas para pagamento (SITUAÇÃO LIQUIDADA)."
$ws.Range("D19").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B20").Value = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D20").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# --- TC3 (rows 23-28) ---
$ws.Range("B25").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B27").Value = "Chefe Clica para exibir a lista de diárias (solicitações) aptas para pagamento (SITUAÇÃO LIQUIDADA)."
$ws.Range("D27").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B28").Value = "Chefe Clica para realizar a autorização de pagamento."
$ws.Range("D28").Value = "SYSTEM Apresenta a tela de Registrar Autorizações de Pagamento."

# --- TC4 (rows 31-36) ---
$ws.Range("B33").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B35").Value = "Chefe Clica para exibir a lista de diárias (solicitações) aptas para pagamento (SITUAÇÃO LIQUIDADA)."
$ws.Range("D35").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B36").Value = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$ws.Range("D36").Value = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

# --- TC5 (rows 39-44) ---
$ws.Range("B41").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B43").Value = "Chefe Clica para exibir a lista de diárias (solicitações) aptas para pagamento (SITUAÇÃO LIQUIDADA)."
$ws.Range("D43").Value = "SYSTEM Recupera e exibe para o usuário a lista de diárias aptas para pagamento ordenado pelo número de diárias em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de autorização (após registrar o empenho)."
$ws.Range("B44").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D44").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
